# Apply a cyclic rotation of data between rows 12, 13 and 14 on the
# "Artfynd" worksheet:
#   new row12 = old row14
#   new row13 = old row12
#   new row14 = old row13
# Only columns A, B, E, F, G, H, P, Q, R, S change (the other columns are
# identical across these three rows, so they are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$cols = @("A", "B", "E", "F", "G", "H", "P", "Q", "R", "S")

# Capture current ("before") values for the three affected rows.
$rows = @(12, 13, 14)
$data = @{}
foreach ($r in $rows) {
    $data[$r] = @{}
    foreach ($c in $cols) {
        $data[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# Source row for each destination row (cyclic shift by one).
$source = @{ 12 = 14; 13 = 12; 14 = 13 }

foreach ($r in $rows) {
    $src = $source[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $data[$src][$c]
    }
}
